$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 141.9375
$ws.Range("I33").Value = 169.33333
$ws.Range("K33").Value = 169.33333
$ws.Range("M33").Value = 59.66667000000001
$ws.Range("H125").Value = 6273.1
$ws.Range("J125").Value = 5979.6
$ws.Range("L125").Value = 53816.4
$ws.Range("N125").Value = -58736.4

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2247.9412
$ws.Range("I32").Value = 2409.4666
$ws.Range("K32").Value = 2409.4666
$ws.Range("M32").Value = -2122.4666
$ws.Range("H132").Value = 439676.1
$ws.Range("I132").Value = 900342.3
$ws.Range("K132").Value = 2701026.9
$ws.Range("M132").Value = -2698496.9
$ws.Range("H139").Value = 58657.316
$ws.Range("J139").Value = 58657.316
$ws.Range("L139").Value = 58657.316
$ws.Range("N139").Value = -68937.31599999999

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 20835162
$ws.Range("I105").Value = 2089.5715
$ws.Range("K105").Value = 2089.5715
$ws.Range("M105").Value = -342.5715
$ws.Range("H134").Value = 45357.69
$ws.Range("I134").Value = 2053.75
$ws.Range("K134").Value = 6161.25
$ws.Range("M134").Value = -3626.25

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 410.3
$ws.Range("J22").Value = 708
$ws.Range("L22").Value = 708
$ws.Range("N22").Value = -1408
$ws.Range("H31").Value = 725004.3
$ws.Range("I31").Value = 1413468.9
$ws.Range("J31").Value = 36539.668
$ws.Range("K31").Value = 1413468.9
$ws.Range("L31").Value = 36539.668
$ws.Range("M31").Value = -1413173.9
$ws.Range("N31").Value = -37129.668
$ws.Range("H34").Value = 725004.3
$ws.Range("I34").Value = 1413468.9
$ws.Range("J34").Value = 36539.668
$ws.Range("K34").Value = 1413468.9
$ws.Range("L34").Value = 36539.668
$ws.Range("M34").Value = -1413266.9
$ws.Range("N34").Value = -36943.668
$ws.Range("H50").Value = 22649.643
$ws.Range("J50").Value = 23001
$ws.Range("L50").Value = 23001
$ws.Range("N50").Value = -24251
$ws.Range("H107").Value = 1347.591
$ws.Range("I107").Value = 883.6667
$ws.Range("K107").Value = 883.6667
$ws.Range("M107").Value = 1036.3333
$ws.Range("H132").Value = 3250.111
$ws.Range("I132").Value = 2039.8948
$ws.Range("K132").Value = 6119.6844
$ws.Range("M132").Value = -3589.6844
$ws.Range("H133").Value = 54024
$ws.Range("J133").Value = 55272.547
$ws.Range("L133").Value = 55272.547
$ws.Range("N133").Value = -60332.547
$ws.Range("H134").Value = 392146.16
$ws.Range("I134").Value = 3046.5293
$ws.Range("K134").Value = 9139.5879
$ws.Range("M134").Value = -6604.5879

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 15989261
$ws.Range("I4").Value = 2263453.8
$ws.Range("K4").Value = 6790361.399999999
$ws.Range("M4").Value = -6790249.399999999
$ws.Range("H5").Value = 1627705.6
$ws.Range("I5").Value = 2760281
$ws.Range("J5").Value = 42100
$ws.Range("K5").Value = 8280843
$ws.Range("L5").Value = 126300
$ws.Range("M5").Value = -8280731
$ws.Range("N5").Value = -126524
$ws.Range("H37").Value = 118900
$ws.Range("J37").Value = 118900
$ws.Range("L37").Value = 356700
$ws.Range("N37").Value = -356924
$ws.Range("H68").Value = 1730.8029
$ws.Range("J68").Value = 1799.4678
$ws.Range("L68").Value = 5398.403399999999
$ws.Range("N68").Value = -7020.403399999999
$ws.Range("H71").Value = 1730.8029
$ws.Range("J71").Value = 1799.4678
$ws.Range("L71").Value = 16195.2102
$ws.Range("N71").Value = -24307.2102
$ws.Range("H75").Value = 1333.3334
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").Value = ""
$ws.Range("H78").Value = 1333.3334
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").Value = ""
$ws.Range("H107").Value = 23693.404
$ws.Range("J107").Value = 25817.209
$ws.Range("L107").Value = 77451.62699999999
$ws.Range("N107").Value = -81291.62699999999
$ws.Range("H122").Value = 846590.5600000001
$ws.Range("I122").Value = 3663370
$ws.Range("K122").Value = 32970330
$ws.Range("M122").Value = -32967880
$ws.Range("H130").Value = 0
$ws.Range("I130").Value = 0
$ws.Range("K130").Value = 0
$ws.Range("M130").Value = ""
$ws.Range("H135").Value = 1627705.6
$ws.Range("I135").Value = 2760281
$ws.Range("J135").Value = 42100
$ws.Range("K135").Value = 24842529
$ws.Range("L135").Value = 378900
$ws.Range("M135").Value = -24839994
$ws.Range("N135").Value = -383970

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 1181763
$ws.Range("I80").Value = 836489.0600000001
$ws.Range("J80").Value = 2010420.4
$ws.Range("K80").Value = 836489.0600000001
$ws.Range("L80").Value = 2010420.4
$ws.Range("M80").Value = -835491.0600000001
$ws.Range("N80").Value = -2012416.4
$ws.Range("H83").Value = 1181763
$ws.Range("I83").Value = 836489.0600000001
$ws.Range("J83").Value = 2010420.4
$ws.Range("K83").Value = 4182445.3
$ws.Range("L83").Value = 10052102
$ws.Range("M83").Value = -4177453.3
$ws.Range("N83").Value = -10062086
$ws.Range("H126").Value = 10066
$ws.Range("J126").Value = 15112.333
$ws.Range("L126").Value = 45336.999
$ws.Range("N126").Value = -50276.999

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2167.25
$ws.Range("I16").Value = 2048.2856
$ws.Range("J16").Value = 3000
$ws.Range("K16").Value = 2048.2856
$ws.Range("L16").Value = 3000
$ws.Range("M16").Value = -1878.2856
$ws.Range("N16").Value = -3340
$ws.Range("H40").Value = 426027.5
$ws.Range("J40").Value = 6747.5
$ws.Range("L40").Value = 6747.5
$ws.Range("N40").Value = -7019.5
$ws.Range("H46").Value = 5032.7334
$ws.Range("I46").Value = 5098.8
$ws.Range("K46").Value = 5098.8
$ws.Range("M46").Value = -4910.8
$ws.Range("H55").Value = 40000640
$ws.Range("J55").Value = 83334330
$ws.Range("L55").Value = 83334330
$ws.Range("N55").Value = -83334676
$ws.Range("H82").Value = 1948.9166
$ws.Range("I82").Value = 950.5
$ws.Range("J82").Value = 2148.6
$ws.Range("K82").Value = 950.5
$ws.Range("L82").Value = 2148.6
$ws.Range("M82").Value = -589.5
$ws.Range("N82").Value = -2870.6
$ws.Range("H85").Value = 1948.9166
$ws.Range("I85").Value = 950.5
$ws.Range("J85").Value = 2148.6
$ws.Range("K85").Value = 950.5
$ws.Range("L85").Value = 2148.6
$ws.Range("M85").Value = 297.5
$ws.Range("N85").Value = -4644.6
$ws.Range("H122").Value = 1381977.1
$ws.Range("I122").Value = 8888
$ws.Range("J122").Value = 1578132.8
$ws.Range("K122").Value = 26664
$ws.Range("L122").Value = 4734398.4
$ws.Range("M122").Value = -24214
$ws.Range("N122").Value = -4739298.4
$ws.Range("H132").Value = 4070.279
$ws.Range("I132").Value = 2927.147
$ws.Range("K132").Value = 8781.440999999999
$ws.Range("M132").Value = -6251.440999999999

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H32").Value = 8332
$ws.Range("I32").Value = 8332
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 8332
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -8015
$ws.Range("N32").Value = ""
$ws.Range("H122").Value = 3595.2666
$ws.Range("I122").Value = 3319
$ws.Range("J122").Value = 4239.8887
$ws.Range("K122").Value = 9957
$ws.Range("L122").Value = 12719.6661
$ws.Range("M122").Value = -7507
$ws.Range("N122").Value = -17619.6661
$ws.Range("H124").Value = 69996.336
$ws.Range("J124").Value = 69996.336
$ws.Range("L124").Value = 69996.336
$ws.Range("N124").Value = -79816.336
$ws.Range("H126").Value = 17844.715
$ws.Range("I126").Value = 5750
$ws.Range("K126").Value = 17250
$ws.Range("M126").Value = -14780
$ws.Range("H132").Value = 28718.068
$ws.Range("I132").Value = 2214.6538
$ws.Range("J132").Value = 67000.78
$ws.Range("K132").Value = 6643.9614
$ws.Range("L132").Value = 201002.34
$ws.Range("M132").Value = -4113.9614
$ws.Range("N132").Value = -206062.34
$ws.Range("H136").Value = 230912.12
$ws.Range("I136").Value = 224534.33
$ws.Range("J136").Value = 266787.25
$ws.Range("K136").Value = 673602.99
$ws.Range("L136").Value = 800361.75
$ws.Range("M136").Value = -801052.36
$ws.Range("N136").Value = -805461.75
